$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Update the summary sheet ("Total") with the new 2022-Q4 row
#    Existing rows 2-7 cascade down by one (B:D only); row 8 is new.
# ---------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Make room for the new row 8 by copying row 7's style into row 8 (A column),
# then cascade the B:D quarter data down one row (bottom-up so nothing is lost).
$total.Range("A7").Copy($total.Range("A8"))
$total.Range("A8").Value = 6

# move the oldest quarter row (2020-Q4, currently row 7) down into the new row 8
$total.Cells.Item(8, 2).Value = $total.Cells.Item(7, 2).Value()
$total.Cells.Item(8, 3).Value = $total.Cells.Item(7, 3).Value()
$total.Cells.Item(8, 4).Value = $total.Cells.Item(7, 4).Value()

# cascade remaining existing quarter rows down by one (keep A column index untouched)
$total.Cells.Item(7, 2).Value = $total.Cells.Item(6, 2).Value()
$total.Cells.Item(7, 3).Value = $total.Cells.Item(6, 3).Value()
$total.Cells.Item(7, 4).Value = $total.Cells.Item(6, 4).Value()
$total.Cells.Item(6, 2).Value = $total.Cells.Item(5, 2).Value()
$total.Cells.Item(6, 3).Value = $total.Cells.Item(5, 3).Value()
$total.Cells.Item(6, 4).Value = $total.Cells.Item(5, 4).Value()
$total.Cells.Item(5, 2).Value = $total.Cells.Item(4, 2).Value()
$total.Cells.Item(5, 3).Value = $total.Cells.Item(4, 3).Value()
$total.Cells.Item(5, 4).Value = $total.Cells.Item(4, 4).Value()
$total.Cells.Item(4, 2).Value = $total.Cells.Item(3, 2).Value()
$total.Cells.Item(4, 3).Value = $total.Cells.Item(3, 3).Value()
$total.Cells.Item(4, 4).Value = $total.Cells.Item(3, 4).Value()
$total.Cells.Item(3, 2).Value = $total.Cells.Item(2, 2).Value()
$total.Cells.Item(3, 3).Value = $total.Cells.Item(2, 3).Value()
$total.Cells.Item(3, 4).Value = $total.Cells.Item(2, 4).Value()

# write the new 2022-Q4 summary row at the top
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 34
$total.Cells.Item(2, 4).Value = 44.91

# ---------------------------------------------------------------
# 2) Insert the new '2022-Q4' detail sheet right before '2022-Q1'
# ---------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")

# Add the new sheet away from '2022-Q1' first (adding it immediately next to
# the sheet it copies styles from loses the copied cell style), fill it in,
# then move it into its final place right before '2022-Q1'.
$q4 = $wb.Worksheets.Add($wb.Worksheets.Item($wb.Worksheets.Count))
$q4.Name = "2022-Q4"

# copy header + index-column styling from the 2022-Q1 sheet so the new sheet
# matches the workbook's look (bold header row, bordered index column)
$q1.Range("A1:H1").Copy($q4.Range("A1:H1"))
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q1.Range("A2").Copy($q4.Range("A2:A35"))

# force the numeric-looking text columns to stay text, matching the source data's
# formatting: fund codes keep leading zeros (B, e.g. '011423') and the scale/
# position/ratio/value columns keep trailing zeros (D:G, e.g. '364.78').
$q4.Range("B2:B35").NumberFormat = "@"
$q4.Range("D2:G35").NumberFormat = "@"

$q4Data = @(
    @(0, '513050', '易方达中证海外中国互联网50（QDII）ETF', '364.78', '98.95', '5.01', '18.2755', 5),
    @(1, '164906', '交银施罗德中证海外中国互联网指数（QDII-LOF）', '113.53', '92.20', '5.19', '5.8922', 5),
    @(2, '159605', '广发中证海外中国互联网30（QDII-ETF）', '58.13', '99.47', '5.83', '3.3890', 6),
    @(3, '011423', '广发全球科技三个月定期开放混合（QDII）美元 C', '25.66', '89.07', '6.99', '1.7936', 1),
    @(4, '011420', '广发全球科技三个月定期开放混合（QDII）人民币 A', '21.02', '89.07', '6.99', '1.4693', 1),
    @(5, '270023', '广发全球精选股票（QDII）', '20.45', '82.63', '6.48', '1.3252', 5),
    @(6, '000906', '广发全球精选股票（QDII）美元现汇', '20.45', '82.63', '6.48', '1.3252', 5),
    @(7, '000934', '国富大中华精选混合（QDII）', '20.61', '87.91', '6.27', '1.2922', 1),
    @(8, '006370', '国富大中华精选混合（QDII）美元', '20.61', '87.91', '6.27', '1.2922', 1),
    @(9, '000988', '嘉实全球互联网股票-人民币（QDII）', '12.08', '89.83', '9.40', '1.1355', 1),
    @(10, '000989', '嘉实全球互联网股票-美元现汇（QDII）', '12.08', '89.83', '9.40', '1.1355', 1),
    @(11, '000990', '嘉实全球互联网股票-美元现钞（QDII）', '12.08', '89.83', '9.40', '1.1355', 1),
    @(12, '159607', '嘉实中证海外中国互联网30ETF（QDII）', '17.02', '99.34', '5.83', '0.9923', 6),
    @(13, '010671', '景顺长城大中华混合（QDII）美元A', '12.44', '86.89', '5.49', '0.6830', 4),
    @(14, '000041', '华夏全球精选股票（QDII）', '18.44', '85.51', '3.70', '0.6823', 4),
    @(15, '262001', '景顺长城大中华混合（QDII）人民币A', '12.42', '86.89', '5.49', '0.6819', 4),
    @(16, '015203', '汇添富全球移动互联灵活配置混合（QDII）D', '11.52', '92.14', '3.45', '0.3974', 6),
    @(17, '001668', '汇添富全球移动互联灵活配置混合（QDII）A', '11.48', '92.14', '3.45', '0.3961', 6),
    @(18, '100055', '富国全球科技互联网股票（QDII）', '3.86', '94.32', '9.14', '0.3528', 1),
    @(19, '011422', '广发全球科技三个月定期开放混合（QDII）人民币 C', '4.84', '89.07', '6.99', '0.3383', 1),
    @(20, '012584', '南方中国新兴经济9个月持有期混合（QDII）A', '3.08', '83.21', '7.17', '0.2208', 1),
    @(21, '457001', '国富亚洲机会股票（QDII）', '3.18', '88.46', '6.24', '0.1984', 2),
    @(22, '006792', '鹏华香港美国互联网股票（LOF）美元现汇', '1.29', '88.46', '7.94', '0.1024', 4),
    @(23, '160644', '鹏华香港美国互联网股票（LOF）人民币', '1.29', '88.46', '7.94', '0.1024', 4),
    @(24, '513220', '招商中证全球中国互联网ETF（QDII）', '1.10', '98.97', '7.13', '0.0784', 5),
    @(25, '012924', '华夏新时代灵活配置混合（QDII）美元现汇', '2.09', '77.57', '2.91', '0.0608', 6),
    @(26, '012925', '华夏新时代灵活配置混合（QDII）美元现钞', '2.09', '77.57', '2.91', '0.0608', 6),
    @(27, '519601', '海富通中国海外精选混合（QDII）', '0.58', '94.54', '6.46', '0.0375', 4),
    @(28, '241001', '华宝海外中国混合（QDII）', '0.76', '93.61', '3.58', '0.0272', 7),
    @(29, '011421', '广发全球科技三个月定期开放混合（QDII）美元 A', '0.20', '89.07', '6.99', '0.0140', 1),
    @(30, '012585', '南方中国新兴经济9个月持有期混合（QDII）C', '0.12', '83.21', '7.17', '0.0086', 1),
    @(31, '519602', '海富通大中华精选混合（QDII）', '0.11', '89.56', '7.52', '0.0083', 1),
    @(32, '016988', '景顺长城大中华混合（QDII）人民币C', '0.02', '86.89', '5.49', '0.0011', 4),
    @(33, '015202', '汇添富全球移动互联灵活配置混合（QDII）C', '0.01', '92.14', '3.45', '0.0003', 6)
)

$r = 2
foreach ($row in $q4Data) {
    $q4.Cells.Item($r, 1).Value = $row[0]
    $q4.Cells.Item($r, 2).Value = $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = $row[3]
    $q4.Cells.Item($r, 5).Value = $row[4]
    $q4.Cells.Item($r, 6).Value = $row[5]
    $q4.Cells.Item($r, 7).Value = $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# finally, relocate the fully-populated sheet to sit right before '2022-Q1'
$q4.Move($q1)
